$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the PREPARATION text (F2) and USERID (G2) with the new userid 31160
# (replacing the previous userid 32382). The FILE_EXCEL cell (O2) keeps its
# text "25012023HargaPasarFixedIncome.xlsx".
$ws.Range("F2").Value = "Username : 31160;`r`nPassword : bni1234;`r`nTgl. Market : 25/01/2023;`r`nFile Excel : 25012023HargaPasarFixedIncome.xlsx"
$ws.Range("G2").Value = 31160

# Move the active selection / view to reflect the saved cursor position.
$ws.Range("F3").Select()
